$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new time-tracking entries (rows 30-32) ---
# Row 30: 28.12.2013, Teresa, 0.5h, "Esoterik - fertig"
$ws.Range("A30").Value = 41636
$ws.Range("B30").Value = "Teresa"
$ws.Range("D30").Value = 0.5
$ws.Range("E30").Value = "Esoterik - fertig"

# Row 31: 02.01.2014, Teresa, 3.5h, "Webserver - post/get"
$ws.Range("A31").Value = 41641
$ws.Range("B31").Value = "Teresa"
$ws.Range("D31").Value = 3.5
$ws.Range("E31").Value = "Webserver - post/get"

# Row 32: 03.01.2014, Teresa, 2h, "Webserver - post/get & Fehler behoben"
$ws.Range("A32").Value = 41642
$ws.Range("B32").Value = "Teresa"
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = "Webserver - post/get & Fehler behoben"

# Copy date formatting from an existing, correctly-formatted date cell (A29)
# onto the newly filled date cells so they pick up the same number format
# (style) rather than creating a brand-new one.
$ws.Range("A29").Copy()
$ws.Range("A30:A32").PasteSpecial(-4122)
$ws.Range("A1").Select()

# --- Update the view so the new rows are visible / selected ---
$ws.Range("C34").Select()
